# Append three new order rows (9, 10, 11) to the sheet, mirroring the
# existing "orderId / customerName / amount / paymentStatus / city" table,
# as part of preparing the data for a MongoDB import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: ORD50013 / narmadha / 2000 / Paid / sangalpet
$ws.Range("A9").Value = "ORD50013"
$ws.Range("B9").Value = "narmadha"
$ws.Range("C9").Value = 2000
$ws.Range("D9").Value = "Paid"
$ws.Range("E9").Value = "sangalpet"

# Row 10: ORD50013 / ashwini / 2500 / Pending / bangalore
$ws.Range("A10").Value = "ORD50013"
$ws.Range("B10").Value = "ashwini"
$ws.Range("C10").Value = 2500
$ws.Range("D10").Value = "Pending"
$ws.Range("E10").Value = "bangalore"

# Row 11: ORD50013 / lokesh / 2500 / Pending / bangalore
$ws.Range("A11").Value = "ORD50013"
$ws.Range("B11").Value = "lokesh"
$ws.Range("C11").Value = 2500
$ws.Range("D11").Value = "Pending"
$ws.Range("E11").Value = "bangalore"

# The sheet originally flags A1:E8 as "number stored as text" (order ids such
# as ORD5006 look numeric-ish) and ignores that warning. Extend the same
# ignored-error range to cover the new rows now that data goes to A1:E11.
try {
    $ws.Range("A1:E11").Errors.Item(6).Ignore = $true
} catch {
    # Older/limited hosts may not expose per-range error suppression; the
    # row data above is the substantive part of this change.
}
